$wb = $excel.ActiveWorkbook

# ============================================================
# Step 1: insert a new worksheet named "2022-Q1" immediately
#         before the "总计" (totals) sheet.
# ============================================================
$totalSheet = $wb.Worksheets.Item('总计')
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = '2022-Q1'

# ============================================================
# Step 2: clone the header-row / index-column cell formatting
#         (bold, centered, bordered style) from the existing
#         "2021-Q4" sheet, which has an identical layout.
# ============================================================
$templateSheet = $wb.Worksheets.Item('2021-Q4')
$templateSheet.Range('B1:H1').Copy()
$newSheet.Range('B1:H1').PasteSpecial(-4122)
$templateSheet.Range('A2').Copy()
$newSheet.Range('A2:A32').PasteSpecial(-4122)

# ============================================================
# Step 3: header row text for "2022-Q1".
# ============================================================
$newSheet.Range('B1').Value = '基金代码'
$newSheet.Range('C1').Value = '基金名称'
$newSheet.Range('D1').Value = '基金规模'
$newSheet.Range('E1').Value = '股票总仓位'
$newSheet.Range('F1').Value = '仓位占比'
$newSheet.Range('G1').Value = '持有市值(亿元)'
$newSheet.Range('H1').Value = '仓位排名'

# ============================================================
# Step 4: fund holdings data rows for "2022-Q1" (A2:H32).
# Columns B and D:G hold text that looks numeric (fund codes,
# percentages truncated to 2 decimals, etc.) in the source data
# -- force them to Text format before assignment so they are not
# auto-converted to numbers (which would also drop leading zeros
# in fund codes such as "000055").
# ============================================================
$newSheet.Range('B2:B32').NumberFormat = '@'
$newSheet.Range('D2:G32').NumberFormat = '@'

$newSheet.Range('A2').Value = 0
$newSheet.Range('B2').Value = '000055'
$newSheet.Range('C2').Value = '广发纳斯达克100指数(QDII) A 美元现汇'
$newSheet.Range('D2').Value = '75.36'
$newSheet.Range('E2').Value = '85.84'
$newSheet.Range('F2').Value = '12.44'
$newSheet.Range('G2').Value = '9.3748'
$newSheet.Range('H2').Value = 1
$newSheet.Range('A3').Value = 1
$newSheet.Range('B3').Value = '270042'
$newSheet.Range('C3').Value = '广发纳斯达克100指数QDII A'
$newSheet.Range('D3').Value = '75.36'
$newSheet.Range('E3').Value = '85.84'
$newSheet.Range('F3').Value = '12.44'
$newSheet.Range('G3').Value = '9.3748'
$newSheet.Range('H3').Value = 1
$newSheet.Range('A4').Value = 2
$newSheet.Range('B4').Value = '006479'
$newSheet.Range('C4').Value = '广发纳斯达克100指数（QDII）C人民币'
$newSheet.Range('D4').Value = '75.36'
$newSheet.Range('E4').Value = '85.84'
$newSheet.Range('F4').Value = '12.44'
$newSheet.Range('G4').Value = '9.3748'
$newSheet.Range('H4').Value = 1
$newSheet.Range('A5').Value = 3
$newSheet.Range('B5').Value = '006480'
$newSheet.Range('C5').Value = '广发纳斯达克100指数（QDII）C美元现汇'
$newSheet.Range('D5').Value = '75.36'
$newSheet.Range('E5').Value = '85.84'
$newSheet.Range('F5').Value = '12.44'
$newSheet.Range('G5').Value = '9.3748'
$newSheet.Range('H5').Value = 1
$newSheet.Range('A6').Value = 4
$newSheet.Range('B6').Value = '513100'
$newSheet.Range('C6').Value = '国泰纳斯达克100 (QDII-ETF)'
$newSheet.Range('D6').Value = '41.86'
$newSheet.Range('E6').Value = '90.74'
$newSheet.Range('F6').Value = '11.42'
$newSheet.Range('G6').Value = '4.7804'
$newSheet.Range('H6').Value = 1
$newSheet.Range('A7').Value = 5
$newSheet.Range('B7').Value = '513500'
$newSheet.Range('C7').Value = '博时标普500ETF(QDII)'
$newSheet.Range('D7').Value = '70.03'
$newSheet.Range('E7').Value = '90.45'
$newSheet.Range('F7').Value = '6.40'
$newSheet.Range('G7').Value = '4.4819'
$newSheet.Range('H7').Value = 1
$newSheet.Range('A8').Value = 6
$newSheet.Range('B8').Value = '040046'
$newSheet.Range('C8').Value = '华安纳斯达克100指数QDII - 人民币'
$newSheet.Range('D8').Value = '22.85'
$newSheet.Range('E8').Value = '90.93'
$newSheet.Range('F8').Value = '11.43'
$newSheet.Range('G8').Value = '2.6118'
$newSheet.Range('H8').Value = 1
$newSheet.Range('A9').Value = 7
$newSheet.Range('B9').Value = '040047'
$newSheet.Range('C9').Value = '华安纳斯达克100指数QDII - 美元现钞'
$newSheet.Range('D9').Value = '22.85'
$newSheet.Range('E9').Value = '90.93'
$newSheet.Range('F9').Value = '11.43'
$newSheet.Range('G9').Value = '2.6118'
$newSheet.Range('H9').Value = 1
$newSheet.Range('A10').Value = 8
$newSheet.Range('B10').Value = '040048'
$newSheet.Range('C10').Value = '华安纳斯达克100指数QDII - 美元现汇'
$newSheet.Range('D10').Value = '22.85'
$newSheet.Range('E10').Value = '90.93'
$newSheet.Range('F10').Value = '11.43'
$newSheet.Range('G10').Value = '2.6118'
$newSheet.Range('H10').Value = 1
$newSheet.Range('A11').Value = 9
$newSheet.Range('B11').Value = '160213'
$newSheet.Range('C11').Value = '国泰纳斯达克100指数(QDII)'
$newSheet.Range('D11').Value = '15.88'
$newSheet.Range('E11').Value = '90.49'
$newSheet.Range('F11').Value = '11.27'
$newSheet.Range('G11').Value = '1.7897'
$newSheet.Range('H11').Value = 1
$newSheet.Range('A12').Value = 10
$newSheet.Range('B12').Value = '000834'
$newSheet.Range('C12').Value = '大成纳斯达克100指数 (QDII)'
$newSheet.Range('D12').Value = '14.15'
$newSheet.Range('E12').Value = '89.27'
$newSheet.Range('F12').Value = '11.17'
$newSheet.Range('G12').Value = '1.5806'
$newSheet.Range('H12').Value = 1
$newSheet.Range('A13').Value = 11
$newSheet.Range('B13').Value = '003721'
$newSheet.Range('C13').Value = '易方达标普信息科技指数（QDII-LOF）美元'
$newSheet.Range('D13').Value = '6.31'
$newSheet.Range('E13').Value = '93.58'
$newSheet.Range('F13').Value = '23.65'
$newSheet.Range('G13').Value = '1.4923'
$newSheet.Range('H13').Value = 1
$newSheet.Range('A14').Value = 12
$newSheet.Range('B14').Value = '161128'
$newSheet.Range('C14').Value = '易方达标普信息科技指数（QDII-LOF）人民币'
$newSheet.Range('D14').Value = '6.31'
$newSheet.Range('E14').Value = '93.58'
$newSheet.Range('F14').Value = '23.65'
$newSheet.Range('G14').Value = '1.4923'
$newSheet.Range('H14').Value = 1
$newSheet.Range('A15').Value = 13
$newSheet.Range('B15').Value = '000043'
$newSheet.Range('C15').Value = '嘉实美国成长股票(QDII) -人民币'
$newSheet.Range('D15').Value = '14.64'
$newSheet.Range('E15').Value = '94.24'
$newSheet.Range('F15').Value = '9.98'
$newSheet.Range('G15').Value = '1.4611'
$newSheet.Range('H15').Value = 1
$newSheet.Range('A16').Value = 14
$newSheet.Range('B16').Value = '000044'
$newSheet.Range('C16').Value = '嘉实美国成长股票(QDII) - 美元现汇'
$newSheet.Range('D16').Value = '14.64'
$newSheet.Range('E16').Value = '94.24'
$newSheet.Range('F16').Value = '9.98'
$newSheet.Range('G16').Value = '1.4611'
$newSheet.Range('H16').Value = 1
$newSheet.Range('A17').Value = 15
$newSheet.Range('B17').Value = '159941'
$newSheet.Range('C17').Value = '广发纳斯达克100ETFQDII'
$newSheet.Range('D17').Value = '11.87'
$newSheet.Range('E17').Value = '90.26'
$newSheet.Range('F17').Value = '12.13'
$newSheet.Range('G17').Value = '1.4398'
$newSheet.Range('H17').Value = 1
$newSheet.Range('A18').Value = 16
$newSheet.Range('B18').Value = '000988'
$newSheet.Range('C18').Value = '嘉实全球互联网股票 - 人民币QDII'
$newSheet.Range('D18').Value = '13.21'
$newSheet.Range('E18').Value = '85.88'
$newSheet.Range('F18').Value = '9.90'
$newSheet.Range('G18').Value = '1.3078'
$newSheet.Range('H18').Value = 4
$newSheet.Range('A19').Value = 17
$newSheet.Range('B19').Value = '000989'
$newSheet.Range('C19').Value = '嘉实全球互联网股票 - 美元现汇QDII'
$newSheet.Range('D19').Value = '13.21'
$newSheet.Range('E19').Value = '85.88'
$newSheet.Range('F19').Value = '9.90'
$newSheet.Range('G19').Value = '1.3078'
$newSheet.Range('H19').Value = 4
$newSheet.Range('A20').Value = 18
$newSheet.Range('B20').Value = '000990'
$newSheet.Range('C20').Value = '嘉实全球互联网股票 - 美元现钞QDII'
$newSheet.Range('D20').Value = '13.21'
$newSheet.Range('E20').Value = '85.88'
$newSheet.Range('F20').Value = '9.90'
$newSheet.Range('G20').Value = '1.3078'
$newSheet.Range('H20').Value = 4
$newSheet.Range('A21').Value = 19
$newSheet.Range('B21').Value = '003722'
$newSheet.Range('C21').Value = '易方达纳斯达克100指数美元（QDII-LOF）'
$newSheet.Range('D21').Value = '9.07'
$newSheet.Range('E21').Value = '91.29'
$newSheet.Range('F21').Value = '11.44'
$newSheet.Range('G21').Value = '1.0376'
$newSheet.Range('H21').Value = 1
$newSheet.Range('A22').Value = 20
$newSheet.Range('B22').Value = '161130'
$newSheet.Range('C22').Value = '易方达纳斯达克100指数人民币（QDII-LOF）'
$newSheet.Range('D22').Value = '9.07'
$newSheet.Range('E22').Value = '91.29'
$newSheet.Range('F22').Value = '11.44'
$newSheet.Range('G22').Value = '1.0376'
$newSheet.Range('H22').Value = 1
$newSheet.Range('A23').Value = 21
$newSheet.Range('B23').Value = '270023'
$newSheet.Range('C23').Value = '广发全球精选股票(QDII)'
$newSheet.Range('D23').Value = '25.53'
$newSheet.Range('E23').Value = '78.43'
$newSheet.Range('F23').Value = '3.58'
$newSheet.Range('G23').Value = '0.9140'
$newSheet.Range('H23').Value = 8
$newSheet.Range('A24').Value = 22
$newSheet.Range('B24').Value = '000906'
$newSheet.Range('C24').Value = '广发全球精选股票(QDII)美元现汇'
$newSheet.Range('D24').Value = '25.53'
$newSheet.Range('E24').Value = '78.43'
$newSheet.Range('F24').Value = '3.58'
$newSheet.Range('G24').Value = '0.9140'
$newSheet.Range('H24').Value = 8
$newSheet.Range('A25').Value = 23
$newSheet.Range('B25').Value = '003718'
$newSheet.Range('C25').Value = '易方达标普500指数(QDII-LOF) 美元'
$newSheet.Range('D25').Value = '5.22'
$newSheet.Range('E25').Value = '91.11'
$newSheet.Range('F25').Value = '6.45'
$newSheet.Range('G25').Value = '0.3367'
$newSheet.Range('H25').Value = 1
$newSheet.Range('A26').Value = 24
$newSheet.Range('B26').Value = '161125'
$newSheet.Range('C26').Value = '易方达标普500指数(QDII-LOF) 人民币'
$newSheet.Range('D26').Value = '5.22'
$newSheet.Range('E26').Value = '91.11'
$newSheet.Range('F26').Value = '6.45'
$newSheet.Range('G26').Value = '0.3367'
$newSheet.Range('H26').Value = 1
$newSheet.Range('A27').Value = 25
$newSheet.Range('B27').Value = '006555'
$newSheet.Range('C27').Value = '浦银安盛全球智能科技股票（QDII）'
$newSheet.Range('D27').Value = '3.20'
$newSheet.Range('E27').Value = '85.41'
$newSheet.Range('F27').Value = '3.46'
$newSheet.Range('G27').Value = '0.1107'
$newSheet.Range('H27').Value = 2
$newSheet.Range('A28').Value = 26
$newSheet.Range('B28').Value = '012924'
$newSheet.Range('C28').Value = '华夏新时代灵活配置混合（QDII）美元现汇'
$newSheet.Range('D28').Value = '2.56'
$newSheet.Range('E28').Value = '84.71'
$newSheet.Range('F28').Value = '3.72'
$newSheet.Range('G28').Value = '0.0952'
$newSheet.Range('H28').Value = 3
$newSheet.Range('A29').Value = 27
$newSheet.Range('B29').Value = '012925'
$newSheet.Range('C29').Value = '华夏新时代灵活配置混合（QDII）美元现钞'
$newSheet.Range('D29').Value = '2.56'
$newSheet.Range('E29').Value = '84.71'
$newSheet.Range('F29').Value = '3.72'
$newSheet.Range('G29').Value = '0.0952'
$newSheet.Range('H29').Value = 3
$newSheet.Range('A30').Value = 28
$newSheet.Range('B30').Value = '005698'
$newSheet.Range('C30').Value = '华夏全球科技先锋混合QDII'
$newSheet.Range('D30').Value = '0.70'
$newSheet.Range('E30').Value = '89.12'
$newSheet.Range('F30').Value = '8.47'
$newSheet.Range('G30').Value = '0.0593'
$newSheet.Range('H30').Value = 4
$newSheet.Range('A31').Value = 29
$newSheet.Range('B31').Value = '519981'
$newSheet.Range('C31').Value = '长信美国标准普尔100等权重指数增强(QDII)'
$newSheet.Range('D31').Value = '0.47'
$newSheet.Range('E31').Value = '84.16'
$newSheet.Range('F31').Value = '0.87'
$newSheet.Range('G31').Value = '0.0041'
$newSheet.Range('H31').Value = 8
$newSheet.Range('A32').Value = 30
$newSheet.Range('B32').Value = '011706'
$newSheet.Range('C32').Value = '长信美国标准普尔100等权重指数增强(QDII) - 美元'
$newSheet.Range('D32').Value = '0.47'
$newSheet.Range('E32').Value = '84.16'
$newSheet.Range('F32').Value = '0.87'
$newSheet.Range('G32').Value = '0.0041'
$newSheet.Range('H32').Value = 8

# Clear the temporary "@" text format back to the default (no explicit
# style index), matching the plain un-styled text cells in the source.
$newSheet.Range('B2:B32').Style = 'Normal'
$newSheet.Range('D2:G32').Style = 'Normal'

# ============================================================
# Step 5: update the "总计" (totals) summary sheet: a new
#         2022-Q1 row is inserted at the top of the data and
#         existing rows shift down with incremented index.
# Re-fetch the sheet by name here: inserting "2022-Q1" shifted
# worksheet positions, and the earlier $totalSheet reference
# resolves by position, not stable identity.
# ============================================================
$totalSheet = $wb.Worksheets.Item('总计')

$totalSheet.Range('A2').Value = 0
$totalSheet.Range('B2').Value = '2022-Q1'
$totalSheet.Range('C2').Value = 31
$totalSheet.Range('D2').Value = 74.18000000000001

$totalSheet.Range('A3').Value = 1
$totalSheet.Range('B3').Value = '2021-Q4'
$totalSheet.Range('C3').Value = 27
$totalSheet.Range('D3').Value = 73.12

$totalSheet.Range('A4').Value = 2
$totalSheet.Range('B4').Value = '2021-Q2'
$totalSheet.Range('C4').Value = 25
$totalSheet.Range('D4').Value = 42.82

$totalSheet.Range('A5').Value = 3
$totalSheet.Range('B5').Value = '2021-Q1'
$totalSheet.Range('C5').Value = 24
$totalSheet.Range('D5').Value = 34.34

$totalSheet.Range('A6').Value = 4
$totalSheet.Range('B6').Value = '2020-Q4'
$totalSheet.Range('C6').Value = 32
$totalSheet.Range('D6').Value = 37.81

# Row 6 is brand new -- clone the index-column style (bold/centered/
# bordered, as used by A2:A5) down onto the new A6 cell.
$totalSheet.Range('A2').Copy()
$totalSheet.Range('A6').PasteSpecial(-4122)
$totalSheet.Range('A6').Value = 4
